$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("O1").Value = "Responsable de facturación"
$ws.Range("P1").Value = "Razón social"

# Match the header style used by the rest of row 1 (e.g. N1 "Tipo de moneda")
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)

# Rows 2-8: column O holds an #N/A error, column P holds "No encontrado"
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 15).Value = "#N/A"
    $ws.Cells.Item($r, 16).Value = "No encontrado"
}

# Rows 9-17: both column O and column P hold "No encontrado"
for ($r = 9; $r -le 17; $r++) {
    $ws.Cells.Item($r, 15).Value = "No encontrado"
    $ws.Cells.Item($r, 16).Value = "No encontrado"
}
